$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.580.54"
$ws.Range("E2").Value = "  +5.65%  "

$ws.Range("D3").Value = "2.739.30"
$ws.Range("E3").Value = "  +4.76%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.38%  "

$ws.Range("E7").Value = "  +1.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.567"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0851"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.19%  "

$ws.Range("E13").Value = "  +2.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.39%  "

$ws.Range("D15").Value = "3.172.29"
$ws.Range("E15").Value = "  +4.88%  "

$ws.Range("D16").Value = "2.769.67"
$ws.Range("E16").Value = "  +6.18%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "51.498.81"
$ws.Range("E17").Value = "  +5.62%  "

$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.872"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.11%  "

$ws.Range("D22").Value = "0.0₃0969"
$ws.Range("E22").Value = "  +2.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("E25").Value = "  +4.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("E34").Value = "  +2.71%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "127.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.72%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0344"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.69%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.81%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.113"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.19%  "

$ws.Range("E45").Value = "  +12.05%  "

$ws.Range("D46").Value = "2.087.13"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.52%  "

$ws.Range("E48").Value = "  +2.88%  "

$ws.Range("E49").Value = "  +6.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.05%  "
